$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Mathew Karl"
$ws.Range("B2").Value = "mathew@nxglabs.in"
$ws.Range("C2").Value = 35534343434

$ws.Range("A3").Value = "Tony Stark"
$ws.Range("B3").Value = "tonys@nxglabs.in"
$ws.Range("C3").Value = 233343434

$ws.Range("A4").Value = "Andy amaya"
$ws.Range("B4").Value = "andyamaya@nxglabs.in"
$ws.Range("C4").Value = 3367546546

$ws.Range("A5").Value = $null
$ws.Range("B5").Value = $null
$ws.Range("C5").Value = $null
$ws.Range("A6").Value = $null
$ws.Range("B6").Value = $null

$ws.Hyperlinks.Delete()

$ws.Range("B7").Select()
